$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

$data = New-Object 'object[,]' 18,20
$data[0,0] = "ECs"
$data[0,1] = "Gnas"
$data[0,2] = "Avpr2"
$data[0,3] = "Inflammatory-Mac"
$data[0,4] = 2
$data[0,5] = 1
$data[0,6] = 236.082283
$data[0,7] = 472.164566
$data[0,8] = 0.356044175747874
$data[0,9] = 0.2923054205394457
$data[0,10] = 2
$data[0,11] = 0.6666666666666666
$data[0,12] = 0.1811646666666666
$data[0,13] = 0.5434939999999999
$data[0,14] = 0.3599820106359796
$data[0,15] = 0.3599820106359796
$data[0,16] = 42.76976810560066
$data[0,17] = 256.618608633604
$data[0,18] = 0.1281694982609498
$data[0,19] = 0.1052246930055852
$data[1,0] = "ECs"
$data[1,1] = "Gnas"
$data[1,2] = "Avpr2"
$data[1,3] = "Neutrophils"
$data[1,4] = 2
$data[1,5] = 1
$data[1,6] = 236.082283
$data[1,7] = 472.164566
$data[1,8] = 0.356044175747874
$data[1,9] = 0.2923054205394457
$data[1,10] = 2
$data[1,11] = 0.6666666666666666
$data[1,12] = 0.234272
$data[1,13] = 0.702816
$data[1,14] = 0.4655085737600355
$data[1,15] = 0.4655085737600355
$data[1,16] = 55.30746860297601
$data[1,17] = 331.844811617856
$data[1,18] = 0.1657416164479602
$data[1,19] = 0.1360706794176448
$data[2,0] = "ECs"
$data[2,1] = "Gnas"
$data[2,2] = "Avpr2"
$data[2,3] = "Resolving-Mac"
$data[2,4] = 2
$data[2,5] = 1
$data[2,6] = 236.082283
$data[2,7] = 472.164566
$data[2,8] = 0.356044175747874
$data[2,9] = 0.2923054205394457
$data[2,10] = 1
$data[2,11] = 0.3333333333333333
$data[2,12] = 0.08782366666666667
$data[2,13] = 0.263471
$data[2,14] = 0.174509415603985
$data[2,15] = 0.174509415603985
$data[2,16] = 20.73361172809767
$data[2,17] = 124.401670368586
$data[2,18] = 0.06213306103896401
$data[2,19] = 0.05101004811621573
$data[3,0] = "FAPs"
$data[3,1] = "Gnas"
$data[3,2] = "Avpr2"
$data[3,3] = "Inflammatory-Mac"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 117.551811
$data[3,7] = 352.655433
$data[3,8] = 0.1772841109604352
$data[3,9] = 0.2183202681257223
$data[3,10] = 2
$data[3,11] = 0.6666666666666666
$data[3,12] = 0.1811646666666666
$data[3,13] = 0.5434939999999999
$data[3,14] = 0.3599820106359796
$data[3,15] = 0.3599820106359796
$data[3,16] = 21.296234655878
$data[3,17] = 191.666111902902
$data[3,18] = 0.06381909071734958
$data[3,19] = 0.0785913690824837
$data[4,0] = "FAPs"
$data[4,1] = "Gnas"
$data[4,2] = "Avpr2"
$data[4,3] = "Neutrophils"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 117.551811
$data[4,7] = 352.655433
$data[4,8] = 0.1772841109604352
$data[4,9] = 0.2183202681257223
$data[4,10] = 2
$data[4,11] = 0.6666666666666666
$data[4,12] = 0.234272
$data[4,13] = 0.702816
$data[4,14] = 0.4655085737600355
$data[4,15] = 0.4655085737600355
$data[4,16] = 27.539097866592
$data[4,17] = 247.851880799328
$data[4,18] = 0.08252727364350806
$data[4,19] = 0.1016299566381135
$data[5,0] = "FAPs"
$data[5,1] = "Gnas"
$data[5,2] = "Avpr2"
$data[5,3] = "Resolving-Mac"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 117.551811
$data[5,7] = 352.655433
$data[5,8] = 0.1772841109604352
$data[5,9] = 0.2183202681257223
$data[5,10] = 1
$data[5,11] = 0.3333333333333333
$data[5,12] = 0.08782366666666667
$data[5,13] = 0.263471
$data[5,14] = 0.174509415603985
$data[5,15] = 0.174509415603985
$data[5,16] = 10.323831065327
$data[5,17] = 92.914479587943
$data[5,18] = 0.03093774659957757
$data[5,19] = 0.03809894240512511
$data[6,0] = "Inflammatory-Mac"
$data[6,1] = "Gnas"
$data[6,2] = "Avpr2"
$data[6,3] = "Inflammatory-Mac"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 65.58120733333334
$data[6,7] = 196.743622
$data[6,8] = 0.09890537575641412
$data[6,9] = 0.1217991168934175
$data[6,10] = 2
$data[6,11] = 0.6666666666666666
$data[6,12] = 0.1811646666666666
$data[6,13] = 0.5434939999999999
$data[6,14] = 0.3599820106359796
$data[6,15] = 0.3599820106359796
$data[6,16] = 11.88099756614089
$data[6,17] = 106.928978095268
$data[6,18] = 0.03560415602750103
$data[6,19] = 0.04384549099297914
$data[7,0] = "Inflammatory-Mac"
$data[7,1] = "Gnas"
$data[7,2] = "Avpr2"
$data[7,3] = "Neutrophils"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 65.58120733333334
$data[7,7] = 196.743622
$data[7,8] = 0.09890537575641412
$data[7,9] = 0.1217991168934175
$data[7,10] = 2
$data[7,11] = 0.6666666666666666
$data[7,12] = 0.234272
$data[7,13] = 0.702816
$data[7,14] = 0.4655085737600355
$data[7,15] = 0.4655085737600355
$data[7,16] = 15.36384060439467
$data[7,17] = 138.274565439552
$data[7,18] = 0.04604130040556872
$data[7,19] = 0.05669853319028662
$data[8,0] = "Inflammatory-Mac"
$data[8,1] = "Gnas"
$data[8,2] = "Avpr2"
$data[8,3] = "Resolving-Mac"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 65.58120733333334
$data[8,7] = 196.743622
$data[8,8] = 0.09890537575641412
$data[8,9] = 0.1217991168934175
$data[8,10] = 1
$data[8,11] = 0.3333333333333333
$data[8,12] = 0.08782366666666667
$data[8,13] = 0.263471
$data[8,14] = 0.174509415603985
$data[8,15] = 0.174509415603985
$data[8,16] = 5.759582092440223
$data[8,17] = 51.83623883196201
$data[8,18] = 0.01725991932334437
$data[8,19] = 0.02125509271015174
$data[9,0] = "MuSCs"
$data[9,1] = "Gnas"
$data[9,2] = "Avpr2"
$data[9,3] = "Inflammatory-Mac"
$data[9,4] = 2
$data[9,5] = 1
$data[9,6] = 137.8159335
$data[9,7] = 275.631867
$data[9,8] = 0.2078451623916705
$data[9,9] = 0.1706368808656166
$data[9,10] = 2
$data[9,11] = 0.6666666666666666
$data[9,12] = 0.1811646666666666
$data[9,13] = 0.5434939999999999
$data[9,14] = 0.3599820106359796
$data[9,15] = 0.3599820106359796
$data[9,16] = 24.96737765388299
$data[9,17] = 149.804265923298
$data[9,18] = 0.07482051945871525
$data[9,19] = 0.06142620746265677
$data[10,0] = "MuSCs"
$data[10,1] = "Gnas"
$data[10,2] = "Avpr2"
$data[10,3] = "Neutrophils"
$data[10,4] = 2
$data[10,5] = 1
$data[10,6] = 137.8159335
$data[10,7] = 275.631867
$data[10,8] = 0.2078451623916705
$data[10,9] = 0.1706368808656166
$data[10,10] = 2
$data[10,11] = 0.6666666666666666
$data[10,12] = 0.234272
$data[10,13] = 0.702816
$data[10,14] = 0.4655085737600355
$data[10,15] = 0.4655085737600355
$data[10,16] = 32.286414372912
$data[10,17] = 193.718486237472
$data[10,18] = 0.0967537051078695
$data[10,19] = 0.07943293104261426
$data[11,0] = "MuSCs"
$data[11,1] = "Gnas"
$data[11,2] = "Avpr2"
$data[11,3] = "Resolving-Mac"
$data[11,4] = 2
$data[11,5] = 1
$data[11,6] = 137.8159335
$data[11,7] = 275.631867
$data[11,8] = 0.2078451623916705
$data[11,9] = 0.1706368808656166
$data[11,10] = 1
$data[11,11] = 0.3333333333333333
$data[11,12] = 0.08782366666666667
$data[11,13] = 0.263471
$data[11,14] = 0.174509415603985
$data[11,15] = 0.174509415603985
$data[11,16] = 12.1035006050595
$data[11,17] = 72.62100363035701
$data[11,18] = 0.03627093782508577
$data[11,19] = 0.02977774236034555
$data[12,0] = "Neutrophils"
$data[12,1] = "Gnas"
$data[12,2] = "Avpr2"
$data[12,3] = "Inflammatory-Mac"
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 39.78016766666667
$data[12,7] = 119.340503
$data[12,8] = 0.05999390055030331
$data[12,9] = 0.07388075774581522
$data[12,10] = 2
$data[12,11] = 0.6666666666666666
$data[12,12] = 0.1811646666666666
$data[12,13] = 0.5434939999999999
$data[12,14] = 0.3599820106359796
$data[12,15] = 0.3599820106359796
$data[12,16] = 7.206760815275778
$data[12,17] = 64.860847337482
$data[12,18] = 0.02159672494599319
$data[12,19] = 0.02659574372064829
$data[13,0] = "Neutrophils"
$data[13,1] = "Gnas"
$data[13,2] = "Avpr2"
$data[13,3] = "Neutrophils"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 39.78016766666667
$data[13,7] = 119.340503
$data[13,8] = 0.05999390055030331
$data[13,9] = 0.07388075774581522
$data[13,10] = 2
$data[13,11] = 0.6666666666666666
$data[13,12] = 0.234272
$data[13,13] = 0.702816
$data[13,14] = 0.4655085737600355
$data[13,15] = 0.4655085737600355
$data[13,16] = 9.319379439605335
$data[13,17] = 83.874414956448
$data[13,18] = 0.0279276750794731
$data[13,19] = 0.03439212616656514
$data[14,0] = "Neutrophils"
$data[14,1] = "Gnas"
$data[14,2] = "Avpr2"
$data[14,3] = "Resolving-Mac"
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 39.78016766666667
$data[14,7] = 119.340503
$data[14,8] = 0.05999390055030331
$data[14,9] = 0.07388075774581522
$data[14,10] = 1
$data[14,11] = 0.3333333333333333
$data[14,12] = 0.08782366666666667
$data[14,13] = 0.263471
$data[14,14] = 0.174509415603985
$data[14,15] = 0.174509415603985
$data[14,16] = 3.493640185101445
$data[14,17] = 31.442761665913
$data[14,18] = 0.01046950052483702
$data[14,19] = 0.0128928878586018
$data[15,0] = "Resolving-Mac"
$data[15,1] = "Gnas"
$data[15,2] = "Avpr2"
$data[15,3] = "Inflammatory-Mac"
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 66.258798
$data[15,7] = 198.776394
$data[15,8] = 0.09992727459330304
$data[15,9] = 0.1230575558299827
$data[15,10] = 2
$data[15,11] = 0.6666666666666666
$data[15,12] = 0.1811646666666666
$data[15,13] = 0.5434939999999999
$data[15,14] = 0.3599820106359796
$data[15,15] = 0.3599820106359796
$data[15,16] = 12.003753053404
$data[15,17] = 108.033777480636
$data[15,18] = 0.03597202122547087
$data[15,19] = 0.04429850637162649
$data[16,0] = "Resolving-Mac"
$data[16,1] = "Gnas"
$data[16,2] = "Avpr2"
$data[16,3] = "Neutrophils"
$data[16,4] = 3
$data[16,5] = 1
$data[16,6] = 66.258798
$data[16,7] = 198.776394
$data[16,8] = 0.09992727459330304
$data[16,9] = 0.1230575558299827
$data[16,10] = 2
$data[16,11] = 0.6666666666666666
$data[16,12] = 0.234272
$data[16,13] = 0.702816
$data[16,14] = 0.4655085737600355
$data[16,15] = 0.4655085737600355
$data[16,16] = 15.522581125056
$data[16,17] = 139.703230125504
$data[16,18] = 0.04651700307565593
$data[16,19] = 0.05728434730481118
$data[17,0] = "Resolving-Mac"
$data[17,1] = "Gnas"
$data[17,2] = "Avpr2"
$data[17,3] = "Resolving-Mac"
$data[17,4] = 3
$data[17,5] = 1
$data[17,6] = 66.258798
$data[17,7] = 198.776394
$data[17,8] = 0.09992727459330304
$data[17,9] = 0.1230575558299827
$data[17,10] = 1
$data[17,11] = 0.3333333333333333
$data[17,12] = 0.08782366666666667
$data[17,13] = 0.263471
$data[17,14] = 0.174509415603985
$data[17,15] = 0.174509415603985
$data[17,16] = 5.819090589286001
$data[17,17] = 52.371815303574
$data[17,18] = 0.01743825029217625
$data[17,19] = 0.02147470215354503

$rng = $ws.Range("A2:T19")
$rng.Value = $data
